# Applies the authored edit described by the commit "Add files via upload":
#   - Arkusz3: rename column header "Liczba" -> "stacje", refresh the city
#     charging-station counts (reorder Szczecin/Poznań/Wrocław, update a few
#     values, and clear the now-unknown counts for the remaining cities).
#   - Arkusz6: relabel the "Moc" (power) categories to the new naming used
#     by the chart (11 kW -> Wallbox, 22 kW -> AC, 50 kW -> DC,
#     350 kW -> Ultraszybka DC).
#   - Selection / active-sheet bookkeeping: the workbook was left with
#     Arkusz6 active (instead of Arkusz3), with a fresh selection on each
#     touched sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Arkusz3 ("Miasto" / "stacje" table)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Arkusz3")

$ws3.Range("B1").Value = "stacje"

$ws3.Range("B3").Value = 119

$ws3.Range("A6").Value = "Szczecin"
$ws3.Range("B6").Value = 51

$ws3.Range("A7").Value = "Poznań"
$ws3.Range("B7").Value = 30

$ws3.Range("B8").ClearContents()
$ws3.Range("B9").ClearContents()
$ws3.Range("B10").ClearContents()

$ws3.Range("A11").Value = "Wrocław"
$ws3.Range("B11").ClearContents()

$ws3.Range("B12").ClearContents()
$ws3.Range("B13").ClearContents()
$ws3.Range("B14").ClearContents()
$ws3.Range("B15").ClearContents()
$ws3.Range("B16").ClearContents()

# ---------------------------------------------------------------------
# Arkusz6 ("Moc" / "czas" table) - relabel charger power categories
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Arkusz6")

$ws6.Range("A2").Value = "Wallbox"
$ws6.Range("A3").Value = "AC"
$ws6.Range("A4").Value = "DC"
$ws6.Range("A5").Value = "Ultraszybka DC"

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 31995
$win.Top = 1410

[void]$ws3.Activate()
[void]$ws3.Range("A1:C17").Select()

[void]$ws6.Activate()
[void]$ws6.Range("C2").Select()
